# Updates cryptos list values (coin name/link/price/volume) to the latest
# scrape snapshot. Row 9 (OKB) drops out of the top list, every row below
# it shifts up by one, and a new row 51 (Cronos) is appended at the bottom;
# several Price/Volume(1h) figures for the untouched rows are refreshed too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D/E hold plain-text figures (e.g. "10.50", "1.980") that Excel
# would otherwise reinterpret/reformat as numbers on assignment, so force
# those specific cells to Text format before writing their new value.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the new cell values cell-by-cell, row by row (A2:A51 / index column
# is untouched by this update).

# Row 2
$ws.Range("D2").Value = "26.372.26"
$ws.Range("E2").Value = "  -1.91%  "

# Row 3
$ws.Range("D3").Value = "1.790.36"
$ws.Range("E3").Value = "  -2.25%  "

# Row 4
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("E5").Value = "  +0.05%  "

# Row 6
$ws.Range("D6").Value = "306.91"
$ws.Range("E6").Value = "  -1.18%  "

# Row 7
$ws.Range("D7").Value = "0.4567"
$ws.Range("E7").Value = "  -0.98%  "

# Row 8
$ws.Range("D8").Value = "0.3633"
$ws.Range("E8").Value = "  -0.56%  "

# Row 9
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.07086"
$ws.Range("E9").Value = "  -1.22%  "

# Row 10
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "0.8755"
$ws.Range("E10").Value = "  -0.35%  "

# Row 11
$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "0.07818"
$ws.Range("E11").Value = "  -0.18%  "

# Row 12
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "19.55"
$ws.Range("E12").Value = "  -0.25%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.798.50"
$ws.Range("E13").Value = "  -2.73%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.273"
$ws.Range("E14").Value = "  -0.90%  "

# Row 15
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "6.323"
$ws.Range("E15").Value = "  -1.01%  "

# Row 16
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "84.92"
$ws.Range("E16").Value = "  -4.37%  "

# Row 17
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").Value = "1.009"
$ws.Range("E17").Value = "  +0.31%  "

# Row 18
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.000008530"
$ws.Range("E18").Value = "  -2.41%  "

# Row 19
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "1.007"
$ws.Range("E19").Value = "  +0.10%  "

# Row 20
$ws.Range("B20").Value = "WrappedBTC"
$ws.Range("C20").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D20").Value = "26.409.60"
$ws.Range("E20").Value = "  -1.87%  "

# Row 21
$ws.Range("D21").Value = "14.24"
$ws.Range("E21").Value = "  -1.83%  "

# Row 22
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "4.984"
$ws.Range("E22").Value = "  -0.23%  "

# Row 23
$ws.Range("B23").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C23").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D23").Value = "2.024.03"
$ws.Range("E23").Value = "  -3.61%  "

# Row 24
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "10.50"
$ws.Range("E24").Value = "  +0.76%  "

# Row 25
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "1.980"
$ws.Range("E25").Value = "  -0.08%  "

# Row 26
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "152.13"
$ws.Range("E26").Value = "  +0.89%  "

# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "17.90"
$ws.Range("E27").Value = "  -1.64%  "

# Row 28
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "2.040"
$ws.Range("E28").Value = "  +2.26%  "

# Row 29
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "112.25"
$ws.Range("E29").Value = "  -1.16%  "

# Row 30
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "4.845"
$ws.Range("E30").Value = "  -2.06%  "

# Row 31
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "0.08671"
$ws.Range("E31").Value = "  -1.58%  "

# Row 32
$ws.Range("B32").Value = "HuobiToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D32").Value = "3.058"
$ws.Range("E32").Value = "  -1.25%  "

# Row 33
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "4.448"
$ws.Range("E33").Value = "  -0.37%  "

# Row 34
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "0.7225"
$ws.Range("E34").Value = "  -5.25%  "

# Row 35
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").Value = "2.663"
$ws.Range("E35").Value = "  +0.41%  "

# Row 36
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "1.110"
$ws.Range("E36").Value = "  -2.51%  "

# Row 37
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "1.005"
$ws.Range("E37").Value = "  -0.04%  "

# Row 38
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "1.076"
$ws.Range("E38").Value = "  -1.31%  "

# Row 39
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.01944"
$ws.Range("E39").Value = "  +0.91%  "

# Row 40
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.05103"
$ws.Range("E40").Value = "  -0.68%  "

# Row 41
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "2.867"
$ws.Range("E41").Value = "  -1.88%  "

# Row 42
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "0.5166"
$ws.Range("E42").Value = "  +3.62%  "

# Row 43
$ws.Range("D43").Value = "6.886"
$ws.Range("E43").Value = "  -0.86%  "

# Row 44
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "0.1516"
$ws.Range("E44").Value = "  -4.96%  "

# Row 45
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "8.025"
$ws.Range("E45").Value = "  -4.08%  "

# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.4692"
$ws.Range("E46").Value = "  +0.49%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "9.883"
$ws.Range("E48").Value = "  -3.61%  "

# Row 49
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "100.07"
$ws.Range("E49").Value = "  -2.36%  "

# Row 50
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "1.589"
$ws.Range("E50").Value = "  -1.25%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.05977"
$ws.Range("E51").Value = "  -2.02%  "
